$d = $word.ActiveDocument

$pairs = @(
    @("2024-07-26 Friday", "2024-07-27 Saturday"),
    @("63÷2=", "68÷6="),
    @("51÷7=", "65÷8="),
    @("44÷8=", "62÷4="),
    @("90÷8=", "99÷7="),
    @("58÷3=", "20÷2="),
    @("26÷9=", "97÷9="),
    @("73÷3=", "23÷9="),
    @("30÷6=", "98÷6="),
    @("45÷4=", "16÷9="),
    @("37÷9=", "66÷3="),
    @("44÷2=", "48÷2="),
    @("89÷7=", "66÷9="),
    @("50÷7=", "69÷3="),
    @("53÷9=", "42÷4="),
    @("57÷9=", "37÷7="),
    @("33÷8=", "97÷9="),
    @("27÷2=", "99÷3="),
    @("13÷5=", "74÷5="),
    @("94÷5=", "65÷3="),
    @("90÷4=", "24÷2="),
    @("80÷3=", "29÷7="),
    @("44÷9=", "46÷5="),
    @("88÷3=", "56÷9="),
    @("98÷5=", "61÷5="),
    @("38÷6=", "41÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
